$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IndividualBiometrics")
$ws.Activate()

# Row 2 ("IndividualId" column) was a leftover placeholder ("Indiv1"); correct it
# to the actual Gender value used for this biometric entry.
$ws.Range("A2").Value = "MALE"

# Add a second individual (a FEMALE counterpart) as a new row.
$ws.Range("A3").Value = "FEMALE"
$ws.Range("B3").Value = "Human"
$ws.Range("C3").Value = "European_ICRP_2002"
$ws.Range("D3").Value = "FEMALE"
$ws.Range("E3").Value = 65
$ws.Range("F3").Value = 165
$ws.Range("G3").Value = 30

[void]$ws.Range("H3").Select()
